$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 125
$ws.Range("C3").Value = 83
$ws.Range("C4").Value = 69
$ws.Range("C5").Value = 43
$ws.Range("C6").Value = 46
$ws.Range("C7").Value = 23
$ws.Range("C8").Value = 29
$ws.Range("C9").Value = 68
$ws.Range("C10").Value = 55
$ws.Range("C11").Value = 32
$ws.Range("C12").Value = 97
$ws.Range("C13").Value = 41
$ws.Range("C14").Value = 71
$ws.Range("C15").Value = 42
$ws.Range("C16").Value = 82
$ws.Range("C17").Value = 126
$ws.Range("C18").Value = 95
$ws.Range("C19").Value = 81
$ws.Range("C20").Value = 53
$ws.Range("C21").Value = 102
$ws.Range("C24").Value = 40
$ws.Range("C25").Value = 28
$ws.Range("C26").Value = 27
$ws.Range("C27").Value = 39
$ws.Range("C29").Value = 57
$ws.Range("C30").Value = 1
$ws.Range("C31").Value = 56
$ws.Range("C32").Value = 113
$ws.Range("C33").Value = 106
$ws.Range("C34").Value = 108
$ws.Range("C35").Value = 103
$ws.Range("C36").Value = 175
$ws.Range("C37").Value = 92
$ws.Range("C38").Value = 100
$ws.Range("C39").Value = 112
$ws.Range("C40").Value = 119
$ws.Range("C41").Value = 65
$ws.Range("C42").Value = 59
